$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: wrap the word "immundices" in editorial <m>...</m> markup (blue,
# Courier New, 9pt) the same way it is already done elsewhere in the
# document for similar words (e.g. "limaille", "ordures").
# ---------------------------------------------------------------------------

# Grab ready-made FormattedText templates for an opening "<m>" tag and a
# closing "</m>" tag from an existing occurrence elsewhere in the document,
# so the new runs get byte-for-byte identical formatting (font, color,
# size) as the rest of the markup tags.
$mOpenSrc = $d.Content.Duplicate
$mOpenSrc.Find.Execute("<m>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mOpenFT = $mOpenSrc.FormattedText

$mCloseSrc = $d.Content.Duplicate
$mCloseSrc.Find.Execute("</m>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mCloseFT = $mCloseSrc.FormattedText

# Locate "immundices" (unique in the document).
$immTarget = $d.Content.Duplicate
$immTarget.Find.Execute("immundices", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$immStart = $immTarget.Start
$immEnd = $immTarget.End

# Insert the closing tag right after "immundices" first, then the opening
# tag right before it, so neither insertion shifts the other's position.
$afterImm = $d.Range($immEnd, $immEnd)
$afterImm.FormattedText = $mCloseFT

$beforeImm = $d.Range($immStart, $immStart)
$beforeImm.FormattedText = $mOpenFT

# ---------------------------------------------------------------------------
# Edit 2: remove the standalone "bonne " run in front of "<m>limaille</m>"
# (the one followed by "ira au fons") and fold its text into the
# "limaille" run so it reads "bonne limaille" instead.
# ---------------------------------------------------------------------------

# Locate "bonne " (unique in the document) and capture its formatting.
$bonneRng = $d.Content.Duplicate
$bonneRng.Find.Execute("bonne ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bonneFT = $bonneRng.FormattedText
$bonneStart = $bonneRng.Start
$bonneEnd = $bonneRng.End

# Find the "limaille" run that immediately follows "bonne " (inside the
# following "<m>...</m>" markup) - search forward from the end of "bonne ".
$limRng = $d.Range($bonneEnd, $d.Content.End)
$limRng.Find.Execute("limaille", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$limStart = $limRng.Start

# Insert "bonne " in front of "limaille" first (it merges with the
# "limaille" run since the formatting is identical), then remove the
# original standalone "bonne " run.
$insPos = $d.Range($limStart, $limStart)
$insPos.FormattedText = $bonneFT

$origBonne = $d.Range($bonneStart, $bonneEnd)
$origBonne.Delete()
